$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new "Wins", "Losses", "Ties" columns (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Data rows 2-61: season record (Wins, Losses, Ties) for every player row
for ($row = 2; $row -le 61; $row++) {
    $ws.Cells.Item($row, 30).Value = 90
    $ws.Cells.Item($row, 31).Value = 72
    $ws.Cells.Item($row, 32).Value = 0
}
